$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130; this shifts existing rows 130-190 down to 131-191,
# carrying forward cell formatting (so the new D130 already inherits the date style).
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with its data.
$ws.Cells.Item(130, 1).Value = 6
$ws.Cells.Item(130, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(130, 3).Value = "Metropolitana"
$ws.Cells.Item(130, 4).Value = (Get-Date -Year 2021 -Month 12 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Cells.Item(130, 5).Value = 13
$ws.Cells.Item(130, 6).Value = 100112022
$ws.Cells.Item(130, 7).Value = "Arveja Verde"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 230
$ws.Cells.Item(130, 11).Value = 10000
$ws.Cells.Item(130, 12).Value = 12000
$ws.Cells.Item(130, 13).Value = 11130
$ws.Cells.Item(130, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(130, 15).Value = "Carahue"
$ws.Cells.Item(130, 16).Value = 445
$ws.Cells.Item(130, 17).Value = 25
$ws.Cells.Item(130, 18).Value = "Hortaliza"
